# RMA Receipt Reversal.xlsx — "RMA Details Maintenance Grid" sheet
# Commit: Added Test plan/Added Custom API-Switch TO Lightning OR classic/
#         Template Changes/WIP-LUI_Work Order
#
# The RMA test-data block (rows 2-4, columns E/F/J) is refreshed from the
# previous "RMA-DFGZ-*" test run to a new "RMA-ZGSG-*" test run (new RMA
# numbers, new sub-RMA numbers, and new Salesforce record ids).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2
$ws.Range("E2").Value = "RMA-ZGSG-001"
$ws.Range("F2").Value = "RMA-ZGSG-1-1"
$ws.Range("J2").Value = "a7s5f000000xNYIAA2"

# Row 3
$ws.Range("E3").Value = "RMA-ZGSG-002"
$ws.Range("F3").Value = "RMA-ZGSG-1-2"
$ws.Range("J3").Value = "a7s5f000000xNYJAA2"

# Row 4
$ws.Range("E4").Value = "RMA-ZGSG-003"
$ws.Range("F4").Value = "RMA-ZGSG-1-3"
$ws.Range("J4").Value = "a7s5f000000xNYKAA2"
